# Update handback-status timestamps for the "Generate Report for Handback" run.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 10:29:00"
$wsZhCn.Range("H2").Value = "2016-03-17 10:29:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 10:29:04"
$wsDeDe.Range("H2").Value = "2016-03-17 10:29:22"
